# Logged Week 16 and performed season sim from Week 17
# - Add new RB player "R.Bonnafon" as a new row on the RB sheet (row 6), stats all zero.
# - Make RB the active sheet/tab (it was WR).
# - Update selections: RB sheet selection -> J7, WR sheet loses tabSelected / selection -> J10.

$wb = $excel.ActiveWorkbook

$rb = $wb.Worksheets.Item("RB")
$wr = $wb.Worksheets.Item("WR")

# Add the new player row to the RB sheet.
$rb.Cells.Item(6, 1).Value = "R.Bonnafon"
$rb.Cells.Item(6, 2).Value = 0
$rb.Cells.Item(6, 3).Value = 0
$rb.Cells.Item(6, 4).Value = 0
$rb.Cells.Item(6, 5).Value = 0
$rb.Cells.Item(6, 6).Value = 0
$rb.Cells.Item(6, 7).Value = 0
$rb.Cells.Item(6, 8).Value = 0
$rb.Cells.Item(6, 9).Value = 0
$rb.Cells.Item(6, 10).Value = 0

# Update the WR sheet's selection before switching away from it.
$wr.Select()
$wr.Range("J10").Select()

# Activate the RB sheet (becomes the workbook's active tab / selected tab).
$rb.Select()
$rb.Range("J7").Select()
